$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Do Before Class" cell for the IPython/Packages/Variables class (row 5).
# Leading apostrophe keeps this a text entry (preserves the existing quote-prefixed
# style already applied to the cell) instead of resetting its formatting.
$ws.Range("D5").Formula = "'- Ipython`n- Packages`n- Python v. R / variables as pointers"

# Update the matching "In-Class Exercise"/prep cell to mention the new Python packages notebook
$ws.Range("E5").Formula = "'- ``Follow this link <https://gke.mybinder.org/v2/gh/ipython/ipython-in-depth/master?filepath=binder/Index.ipynb>``_ , then click `"Ipython - Beyond plain python`" and read that notebook. `n- ``Python packages <managing_python_packages.ipynb>``_`n- ``variables v objects <python_v_r.ipynb>``_"

# Row height grows to fit the extra line of text
$ws.Rows.Item(5).RowHeight = 85

# Update the saved view state (selection moved, no frozen/scrolled top-left row anymore)
$ws.Range("E4").Select()
